# Update for Tidsregistrering PM10 Nikolaj
# Fills in three new time-registration entries (rows 29-31) on the active
# sheet ("Ark1") and updates the current selection/active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Row 29: "Lavet implement af OC0803" / "implenter" -----------------
$ws.Range("A29").Value = "Lavet implement af OC0803"
$ws.Range("B29").Value = "implenter"
$ws.Range("C29").Value = 43896
$ws.Range("D29").Value = 0.375
$ws.Range("E29").Value = 0.47222222222222227

# --- Row 30: "Lavet UI UC04 " / "Designer" ------------------------------
$ws.Range("A30").Value = "Lavet UI UC04 "
$ws.Range("B30").Value = "Designer"
$ws.Range("C30").Value = 43896
$ws.Range("D30").Value = 0.52083333333333337
$ws.Range("E30").Value = 0.5625

# --- Row 31: "Lavet implement af OC0803" / "implenter" ------------------
$ws.Range("A31").Value = "Lavet implement af OC0803"
$ws.Range("B31").Value = "implenter"
$ws.Range("C31").Value = 43896
$ws.Range("D31").Value = 0.5625
$ws.Range("E31").Value = 0.61458333333333337

# --- Update the view state to match where the author ended up working ---
$ws.Range("C32").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 2
